$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 values
$ws.Range("B2").Value = 22.384284072851585
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 23.381542284777765
$ws.Range("E2").Value = 23.260940587398579

# Row 3 values
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = 21.188901479510321
$ws.Range("D3").Value = 14.184921720323047
$ws.Range("E3").Value = 24.188139510940232

# Update selection to match new used range of interest
$ws.Range("B1:E3").Select()
